$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = 0.8077876666666667
$ws.Range("N2").Value = 2.423363
$ws.Range("O2").Value = 0.04902038147436601
$ws.Range("P2").Value = 0.04902038147436601
$ws.Range("Q2").Value = 8.428653075665556
$ws.Range("R2").Value = 75.85787768099
$ws.Range("S2").Value = 0.04760414333330192
$ws.Range("T2").Value = 0.04760414333330194

$ws.Range("O3").Value = 0.4722854529078861
$ws.Range("P3").Value = 0.4722854529078861
$ws.Range("S3").Value = 0.4586407473433722
$ws.Range("T3").Value = 0.4586407473433723

$ws.Range("M4").Value = 5.009781333333333
$ws.Range("N4").Value = 15.029344
$ws.Range("O4").Value = 0.3040172587389813
$ws.Range("P4").Value = 0.3040172587389813
$ws.Range("Q4").Value = 52.2732774787911
$ws.Range("R4").Value = 470.4594973091199
$ws.Range("S4").Value = 0.2952339562754326
$ws.Range("T4").Value = 0.2952339562754326

$ws.Range("M5").Value = 2.878432333333333
$ws.Range("N5").Value = 8.635297
$ws.Range("O5").Value = 0.1746769068787666
$ws.Range("P5").Value = 0.1746769068787666
$ws.Range("Q5").Value = 30.03426338453444
$ws.Range("R5").Value = 270.30837046081
$ws.Range("S5").Value = 0.1696303509270514
$ws.Range("T5").Value = 0.1696303509270514

$ws.Range("M6").Value = 0.8077876666666667
$ws.Range("N6").Value = 2.423363
$ws.Range("O6").Value = 0.04902038147436601
$ws.Range("P6").Value = 0.04902038147436601
$ws.Range("Q6").Value = 0.250755063062
$ws.Range("R6").Value = 2.256795567558
$ws.Range("S6").Value = 0.001416238141064078
$ws.Range("T6").Value = 0.001416238141064078

$ws.Range("O7").Value = 0.4722854529078861
$ws.Range("P7").Value = 0.4722854529078861
$ws.Range("S7").Value = 0.0136447055645138
$ws.Range("T7").Value = 0.0136447055645138

$ws.Range("M8").Value = 5.009781333333333
$ws.Range("N8").Value = 15.029344
$ws.Range("O8").Value = 0.3040172587389813
$ws.Range("P8").Value = 0.3040172587389813
$ws.Range("Q8").Value = 1.555146341056
$ws.Range("R8").Value = 13.996317069504
$ws.Range("S8").Value = 0.008783302463548612
$ws.Range("T8").Value = 0.008783302463548612

$ws.Range("M9").Value = 2.878432333333333
$ws.Range("N9").Value = 8.635297
$ws.Range("O9").Value = 0.1746769068787666
$ws.Range("P9").Value = 0.1746769068787666
$ws.Range("Q9").Value = 0.8935287217779998
$ws.Range("R9").Value = 8.041758496001998
$ws.Range("S9").Value = 0.005046555951715121
$ws.Range("T9").Value = 0.005046555951715121

